# Instalador de base de datos terminado
# Se termina el instalador de base de datos y se implementa la eliminacion
# de un registro de renta.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Fila 12 - "Instalador de aplicacion para Linux.": se registran 4 horas
# consumidas el dia 5 (columna T).
$ws.Range("T12").Value = 4

# Fila 17 - "Prueba general de sistema en Linux.": se marca como Hecho y se
# registran 3 horas consumidas el dia 5 (columna T).
$ws.Range("F17").Value = "Hecho"
$ws.Range("T17").Value = 3

# Se vuelven a fusionar las celdas de encabezado de los ultimos dias
# (efecto colateral de los ajustes realizados sobre esa zona de la tabla).
$ranges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($r in $ranges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $ranges) {
    $ws.Range($r).Merge()
}

# Selecciona la celda activa final tal como quedo al guardar el archivo
$ws.Activate()
$ws.Range("T15").Select()
